$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is being inserted after the existing row 114
# (old row 114 -> "Precio" updated to the new week's figures; a fresh row
# is opened at 115 holding the values the old row 114 used to carry; the
# row that used to be 115 shifts down to 116 unchanged).
$ws.Rows("115").Insert()

# Row 115 (newly opened) takes on the values that row 114 held before
# this edit.
$ws.Range("A115").Value = 10
$ws.Range("B115").Value = "Vega Modelo de Temuco"
$ws.Range("C115").Value = "La Araucanía"
$ws.Range("D115").Value = 44496
$ws.Range("E115").Value = 9
$ws.Range("F115").Value = 100112042
$ws.Range("G115").Value = "Locoto"
$ws.Range("H115").Value = "Sin especificar"
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 40
$ws.Range("K115").Value = 2200
$ws.Range("L115").Value = 2200
$ws.Range("M115").Value = 2200
$ws.Range("N115").Value = "$/kilo"
$ws.Range("O115").Value = "Región de Arica y Parinacota"
$ws.Range("P115").Value = 2200
$ws.Range("Q115").Value = 1
$ws.Range("R115").Value = "Hortaliza"

# Row 114 is updated in place with the new week's price/volume figures.
$ws.Range("D114").Value = 45267
$ws.Range("J114").Value = 100
$ws.Range("K114").Value = 2700
$ws.Range("L114").Value = 2700
$ws.Range("M114").Value = 2700
$ws.Range("P114").Value = 2700
